$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G ("Unit Of Measure*" is F, so the new
# "Other Unit Of Measure" column lands at G) - this shifts the old G..L
# (Weight Per Item*..DangerousGoodsCode) right to H..M, carrying their data
# and column-width formatting with them automatically.
$ws.Columns("G").Insert()

# New header cell for the inserted column.
$ws.Range("G1").Value = "Other Unit Of Measure"

# Best-effort column width for the new column (engine rounds ColumnWidth to
# 1/6 character-width increments, so this lands close to the authored
# 21.5703125 stored width).
$ws.Range("G1").ColumnWidth = 20.7

# New "Point Reference*" values for the two data rows.
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 2

# Selection moved to F3 in the saved file.
$ws.Range("F3").Select()

# Data validation layout changed: the "Unit Of Measure*" list validation
# used to cover F2:F4; now F2/F3 (already filled in) just allow any value,
# while the list dropdown follows into the new column for rows 3:4, and G2
# is added to the "any value" group alongside C2:D5.
$ws.Range("F2:F3").Validation.Delete()
$ws.Range("F2:F3").Validation.Add(0, 1, 1)
$ws.Range("G2").Validation.Add(0, 1, 1)
$ws.Range("G3:G4").Validation.Add(3, 1, 1, """Litre, Box, Bag, Piece, weight-kg""")
